# Apply capital-structure database update for Russia Auto & Truck sheet
# (row 2 metric refresh; rows 3/4 swapped between Sollers (SVAV) and ZIL (ZILL)
# with refreshed metrics; stale "expected_growth_eps_next_5_years" (col F) cleared)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.159
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = 0.04727629588764241
$ws.Range("H2").Value = 0.03649095229509736
$ws.Range("I2").Value = -0.02446405263334043
$ws.Range("J2").Value = -0.02446405263334043
$ws.Range("K2").Value = -19.34
$ws.Range("L2").Value = -0.02087973424222435
$ws.Range("U2").Value = 106.9
$ws.Range("V2").Value = 0.4962859795728877
$ws.Range("W2").Value = 0.3184989736619852
$ws.Range("X2").Value = 0.1751585896915784
$ws.Range("Y2").Value = 0.1433403839704068
$ws.Range("Z2").Value = 1.428439022885695
$ws.Range("AA2").Value = -0.03338947895082806
$ws.Range("AB2").Value = 0.08639540415900573
$ws.Range("AC2").Value = -0.1197848831098338
$ws.Range("AD2").Value = 359.1
$ws.Range("AF2").Value = 359.1
$ws.Range("AG2").Value = 252.2
$ws.Range("AH2").Value = 0.6250652741514361
$ws.Range("AI2").Value = 0.4995270420654351
$ws.Range("AJ2").Value = 0.5393498716852011
$ws.Range("AK2").Value = 0.4121049707506782
$ws.Range("AL2").Value = 36.93
$ws.Range("AM2").Value = 35.02
$ws.Range("AN2").Value = 21.04923798358734
$ws.Range("AO2").Value = -0.6135932845924722
$ws.Range("AP2").Value = 14.7831184056272
$ws.Range("AQ2").Value = -0.6470588235294117
# Row 3
$ws.Range("B3").Value = 'Public Joint-Stock Company The Likhachov Plant (MISX:ZILL)'
$ws.Range("D3").Value = -0.423
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = -4.784130688448075
$ws.Range("H3").Value = -4.784130688448075
$ws.Range("I3").Value = -3.103850641773629
$ws.Range("J3").Value = -3.103850641773629
$ws.Range("K3").Value = -6.64
$ws.Range("L3").Value = -7.747957992998833
$ws.Range("U3").Value = 27.3
$ws.Range("V3").Value = 0.2948164146868251
$ws.Range("W3").Value = 0.680327868852459
$ws.Range("X3").Value = 0.2050379392813463
$ws.Range("Y3").Value = 0.4752899295711127
$ws.Range("Z3").Value = 0.002821492065582406
$ws.Range("AA3").Value = -0.008757489958517153
$ws.Range("AB3").Value = 0.08997256237260873
$ws.Range("AC3").Value = -0.09873005233112588
$ws.Range("AD3").Value = 211.3
$ws.Range("AF3").Value = 211.3
$ws.Range("AG3").Value = 184
$ws.Range("AH3").Value = 0.6952945047713064
$ws.Range("AI3").Value = 1.020871581795342
$ws.Range("AJ3").Value = 0.6652205350686912
$ws.Range("AK3").Value = 1.024042742653606
$ws.Range("AL3").Value = 8.33
$ws.Range("AM3").Value = 6.42
$ws.Range("AN3").Value = -80.03787878787878
$ws.Range("AO3").Value = -0.319327731092437
$ws.Range("AP3").Value = -69.69696969696969
$ws.Range("AQ3").Value = -0.4143302180685359
# Row 4
$ws.Range("B4").Value = 'Sollers Public Joint Stock Company (MISX:SVAV)'
$ws.Range("D4").Value = 0.105
$ws.Range("G4").Value = 0.05175059433758375
$ws.Range("H4").Value = 0.04095526258915064
$ws.Range("I4").Value = -0.02161227577263886
$ws.Range("J4").Value = -0.02161227577263886
$ws.Range("K4").Value = -12.7
$ws.Range("L4").Value = -0.01372379511562567
$ws.Range("U4").Value = 79.59999999999999
$ws.Range("V4").Value = 0.6482084690553745
$ws.Range("W4").Value = -0.04332992152848857
$ws.Range("X4").Value = 0.1452792401018105
$ws.Range("Y4").Value = -0.1886091616302991
$ws.Range("Z4").Value = 2.684653321729039
$ws.Range("AA4").Value = -0.05802146794313896
$ws.Range("AB4").Value = 0.08281824594540273
$ws.Range("AC4").Value = -0.1408397138885417
$ws.Range("AD4").Value = 147.8
$ws.Range("AF4").Value = 147.8
$ws.Range("AG4").Value = 68.20000000000002
$ws.Range("AH4").Value = 0.5461936437546193
$ws.Range("AI4").Value = 0.2887282672396952
$ws.Range("AJ4").Value = 0.3570680628272252
$ws.Range("AK4").Value = 0.1577608142493639
$ws.Range("AL4").Value = 28.6
$ws.Range("AM4").Value = 28.6
$ws.Range("AN4").Value = 7.502538071065991
$ws.Range("AO4").Value = -0.6993006993006993
$ws.Range("AP4").Value = 3.461928934010153
$ws.Range("AQ4").Value = -0.6993006993006993
